$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header C12 from "brojRadnika" to "kapacitetRadnika"
$ws.Range("C12").Value = "kapacitetRadnika"

# New request fields: agencija (I column) can send a "ponuda" (offer) to admin,
# klijent (C column) can send a "zahtev" (request) to agencija for cooperation
$ws.Range("I9").Value = "ponuda"
$ws.Range("C13").Value = "zahtev"

# Column C got a bit wider to fit the new, longer labels
$ws.Columns.Item(3).ColumnWidth = 17.6666666666667

# Selection moved to D5
$ws.Range("D5").Select()
